$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.931.64"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").Value = "1.677.44"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9973"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3655"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.79%  "

$ws.Range("D15").Value = "1.676.41"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.641"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001054"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06599"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9978"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.23"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.28%  "

$ws.Range("D24").Value = "24.954.43"
$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.460"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.423"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.11%  "

$ws.Range("D29").Value = "1.862.99"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.190"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.081"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08468"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.202"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06071"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("E40").Value = "  +2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2094"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9976"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5971"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.847"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5731"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.969"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07034"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.14%  "
